# Auto-generated script to append new sensor log rows to 4 worksheets
# matching the target diff (SeniorConnect_MasterLog.xlsx update).
$wb = $excel.ActiveWorkbook

# --- Sheet: mmWave(InBed) (Worksheets.Item(10)) ---
$ws = $wb.Worksheets.Item(10)
# row 87
$ws.Cells.Item(87, 1).NumberFormat = '@'
$ws.Cells.Item(87, 1).Value = '2026-02-01'
$ws.Cells.Item(87, 2).Value = '21:05:30'
$ws.Cells.Item(87, 3).Value = '21:00'
$ws.Cells.Item(87, 4).Value = 'Bedroom'
$ws.Cells.Item(87, 5).Value = 'In Bed'
$ws.Cells.Item(87, 6).Value = 'Occupied'
# row 88
$ws.Cells.Item(88, 1).NumberFormat = '@'
$ws.Cells.Item(88, 1).Value = '2026-02-01'
$ws.Cells.Item(88, 2).Value = '21:05:31'
$ws.Cells.Item(88, 3).Value = '21:00'
$ws.Cells.Item(88, 4).Value = 'Bedroom'
$ws.Cells.Item(88, 5).Value = 'In Bed'
$ws.Cells.Item(88, 6).Value = 'Occupied'
# row 89
$ws.Cells.Item(89, 1).NumberFormat = '@'
$ws.Cells.Item(89, 1).Value = '2026-02-01'
$ws.Cells.Item(89, 2).Value = '21:05:36'
$ws.Cells.Item(89, 3).Value = '21:00'
$ws.Cells.Item(89, 4).Value = 'Bedroom'
$ws.Cells.Item(89, 5).Value = 'In Bed'
$ws.Cells.Item(89, 6).Value = 'Occupied'
# row 90
$ws.Cells.Item(90, 1).NumberFormat = '@'
$ws.Cells.Item(90, 1).Value = '2026-02-01'
$ws.Cells.Item(90, 2).Value = '21:05:48'
$ws.Cells.Item(90, 3).Value = '21:00'
$ws.Cells.Item(90, 4).Value = 'Bedroom'
$ws.Cells.Item(90, 5).Value = 'In Bed'
$ws.Cells.Item(90, 6).Value = 'Occupied'
# row 91
$ws.Cells.Item(91, 1).NumberFormat = '@'
$ws.Cells.Item(91, 1).Value = '2026-02-01'
$ws.Cells.Item(91, 2).Value = '21:05:52'
$ws.Cells.Item(91, 3).Value = '21:00'
$ws.Cells.Item(91, 4).Value = 'Bedroom'
$ws.Cells.Item(91, 5).Value = 'In Bed'
$ws.Cells.Item(91, 6).Value = 'Occupied'
# row 92
$ws.Cells.Item(92, 1).NumberFormat = '@'
$ws.Cells.Item(92, 1).Value = '2026-02-01'
$ws.Cells.Item(92, 2).Value = '21:06:01'
$ws.Cells.Item(92, 3).Value = '21:00'
$ws.Cells.Item(92, 4).Value = 'Bedroom'
$ws.Cells.Item(92, 5).Value = 'In Bed'
$ws.Cells.Item(92, 6).Value = 'Occupied'
# row 93
$ws.Cells.Item(93, 1).NumberFormat = '@'
$ws.Cells.Item(93, 1).Value = '2026-02-01'
$ws.Cells.Item(93, 2).Value = '21:06:02'
$ws.Cells.Item(93, 3).Value = '21:00'
$ws.Cells.Item(93, 4).Value = 'Bedroom'
$ws.Cells.Item(93, 5).Value = 'In Bed'
$ws.Cells.Item(93, 6).Value = 'Occupied'
# row 94
$ws.Cells.Item(94, 1).NumberFormat = '@'
$ws.Cells.Item(94, 1).Value = '2026-02-01'
$ws.Cells.Item(94, 2).Value = '21:06:03'
$ws.Cells.Item(94, 3).Value = '21:00'
$ws.Cells.Item(94, 4).Value = 'Bedroom'
$ws.Cells.Item(94, 5).Value = 'In Bed'
$ws.Cells.Item(94, 6).Value = 'Occupied'
# row 95
$ws.Cells.Item(95, 1).NumberFormat = '@'
$ws.Cells.Item(95, 1).Value = '2026-02-01'
$ws.Cells.Item(95, 2).Value = '21:06:07'
$ws.Cells.Item(95, 3).Value = '21:00'
$ws.Cells.Item(95, 4).Value = 'Bedroom'
$ws.Cells.Item(95, 5).Value = 'In Bed'
$ws.Cells.Item(95, 6).Value = 'Occupied'
# row 96
$ws.Cells.Item(96, 1).NumberFormat = '@'
$ws.Cells.Item(96, 1).Value = '2026-02-01'
$ws.Cells.Item(96, 2).Value = '21:06:11'
$ws.Cells.Item(96, 3).Value = '21:00'
$ws.Cells.Item(96, 4).Value = 'Bedroom'
$ws.Cells.Item(96, 5).Value = 'In Bed'
$ws.Cells.Item(96, 6).Value = 'Occupied'
# row 97
$ws.Cells.Item(97, 1).NumberFormat = '@'
$ws.Cells.Item(97, 1).Value = '2026-02-01'
$ws.Cells.Item(97, 2).Value = '21:06:13'
$ws.Cells.Item(97, 3).Value = '21:00'
$ws.Cells.Item(97, 4).Value = 'Bedroom'
$ws.Cells.Item(97, 5).Value = 'In Bed'
$ws.Cells.Item(97, 6).Value = 'Occupied'
# row 98
$ws.Cells.Item(98, 1).NumberFormat = '@'
$ws.Cells.Item(98, 1).Value = '2026-02-01'
$ws.Cells.Item(98, 2).Value = '21:06:14'
$ws.Cells.Item(98, 3).Value = '21:00'
$ws.Cells.Item(98, 4).Value = 'Bedroom'
$ws.Cells.Item(98, 5).Value = 'In Bed'
$ws.Cells.Item(98, 6).Value = 'Occupied'
# row 99
$ws.Cells.Item(99, 1).NumberFormat = '@'
$ws.Cells.Item(99, 1).Value = '2026-02-01'
$ws.Cells.Item(99, 2).Value = '21:06:24'
$ws.Cells.Item(99, 3).Value = '21:00'
$ws.Cells.Item(99, 4).Value = 'Bedroom'
$ws.Cells.Item(99, 5).Value = 'In Bed'
$ws.Cells.Item(99, 6).Value = 'Occupied'

# --- Sheet: Proximity (Worksheets.Item(5)) ---
$ws = $wb.Worksheets.Item(5)
# row 7
$ws.Cells.Item(7, 1).NumberFormat = '@'
$ws.Cells.Item(7, 1).Value = '2026-02-01'
$ws.Cells.Item(7, 2).Value = '21:05:30'
$ws.Cells.Item(7, 3).Value = '21:00'
$ws.Cells.Item(7, 4).Value = 'Bedroom Door'
$ws.Cells.Item(7, 5).Value = 'ENTER'
$ws.Cells.Item(7, 6).Value = 'User ENTERED Bedroom'
# row 8
$ws.Cells.Item(8, 1).NumberFormat = '@'
$ws.Cells.Item(8, 1).Value = '2026-02-01'
$ws.Cells.Item(8, 2).Value = '21:05:32'
$ws.Cells.Item(8, 3).Value = '21:00'
$ws.Cells.Item(8, 4).Value = 'Bedroom Door'
$ws.Cells.Item(8, 5).Value = 'EXIT'
$ws.Cells.Item(8, 6).Value = 'User EXITED Bedroom'
# row 9
$ws.Cells.Item(9, 1).NumberFormat = '@'
$ws.Cells.Item(9, 1).Value = '2026-02-01'
$ws.Cells.Item(9, 2).Value = '21:05:36'
$ws.Cells.Item(9, 3).Value = '21:00'
$ws.Cells.Item(9, 4).Value = 'Bedroom Door'
$ws.Cells.Item(9, 5).Value = 'ENTER'
$ws.Cells.Item(9, 6).Value = 'User ENTERED Bedroom'

# --- Sheet: mmWave(BR) (Worksheets.Item(8)) ---
$ws = $wb.Worksheets.Item(8)
# row 82
$ws.Cells.Item(82, 1).NumberFormat = '@'
$ws.Cells.Item(82, 1).Value = '2026-02-01'
$ws.Cells.Item(82, 2).Value = '21:05:29'
$ws.Cells.Item(82, 3).Value = '21:00'
$ws.Cells.Item(82, 4).Value = 'Bedroom'
$ws.Cells.Item(82, 5).Value = 12
$ws.Cells.Item(82, 6).Value = 'Occupied'
# row 83
$ws.Cells.Item(83, 1).NumberFormat = '@'
$ws.Cells.Item(83, 1).Value = '2026-02-01'
$ws.Cells.Item(83, 2).Value = '21:05:31'
$ws.Cells.Item(83, 3).Value = '21:00'
$ws.Cells.Item(83, 4).Value = 'Bedroom'
$ws.Cells.Item(83, 5).Value = 5
$ws.Cells.Item(83, 6).Value = 'Occupied'
# row 84
$ws.Cells.Item(84, 1).NumberFormat = '@'
$ws.Cells.Item(84, 1).Value = '2026-02-01'
$ws.Cells.Item(84, 2).Value = '21:05:31'
$ws.Cells.Item(84, 3).Value = '21:00'
$ws.Cells.Item(84, 4).Value = 'Bedroom'
$ws.Cells.Item(84, 5).Value = 2
$ws.Cells.Item(84, 6).Value = 'Occupied'
# row 85
$ws.Cells.Item(85, 1).NumberFormat = '@'
$ws.Cells.Item(85, 1).Value = '2026-02-01'
$ws.Cells.Item(85, 2).Value = '21:05:37'
$ws.Cells.Item(85, 3).Value = '21:00'
$ws.Cells.Item(85, 4).Value = 'Bedroom'
$ws.Cells.Item(85, 5).Value = 1
$ws.Cells.Item(85, 6).Value = 'Occupied'
# row 86
$ws.Cells.Item(86, 1).NumberFormat = '@'
$ws.Cells.Item(86, 1).Value = '2026-02-01'
$ws.Cells.Item(86, 2).Value = '21:05:49'
$ws.Cells.Item(86, 3).Value = '21:00'
$ws.Cells.Item(86, 4).Value = 'Bedroom'
$ws.Cells.Item(86, 5).Value = 2
$ws.Cells.Item(86, 6).Value = 'Occupied'
# row 87
$ws.Cells.Item(87, 1).NumberFormat = '@'
$ws.Cells.Item(87, 1).Value = '2026-02-01'
$ws.Cells.Item(87, 2).Value = '21:05:53'
$ws.Cells.Item(87, 3).Value = '21:00'
$ws.Cells.Item(87, 4).Value = 'Bedroom'
$ws.Cells.Item(87, 5).Value = 1
$ws.Cells.Item(87, 6).Value = 'Occupied'
# row 88
$ws.Cells.Item(88, 1).NumberFormat = '@'
$ws.Cells.Item(88, 1).Value = '2026-02-01'
$ws.Cells.Item(88, 2).Value = '21:06:02'
$ws.Cells.Item(88, 3).Value = '21:00'
$ws.Cells.Item(88, 4).Value = 'Bedroom'
$ws.Cells.Item(88, 5).Value = 19
$ws.Cells.Item(88, 6).Value = 'Occupied'
# row 89
$ws.Cells.Item(89, 1).NumberFormat = '@'
$ws.Cells.Item(89, 1).Value = '2026-02-01'
$ws.Cells.Item(89, 2).Value = '21:06:03'
$ws.Cells.Item(89, 3).Value = '21:00'
$ws.Cells.Item(89, 4).Value = 'Bedroom'
$ws.Cells.Item(89, 5).Value = 3
$ws.Cells.Item(89, 6).Value = 'Occupied'
# row 90
$ws.Cells.Item(90, 1).NumberFormat = '@'
$ws.Cells.Item(90, 1).Value = '2026-02-01'
$ws.Cells.Item(90, 2).Value = '21:06:04'
$ws.Cells.Item(90, 3).Value = '21:00'
$ws.Cells.Item(90, 4).Value = 'Bedroom'
$ws.Cells.Item(90, 5).Value = 2
$ws.Cells.Item(90, 6).Value = 'Occupied'
# row 91
$ws.Cells.Item(91, 1).NumberFormat = '@'
$ws.Cells.Item(91, 1).Value = '2026-02-01'
$ws.Cells.Item(91, 2).Value = '21:06:08'
$ws.Cells.Item(91, 3).Value = '21:00'
$ws.Cells.Item(91, 4).Value = 'Bedroom'
$ws.Cells.Item(91, 5).Value = 1
$ws.Cells.Item(91, 6).Value = 'Occupied'
# row 92
$ws.Cells.Item(92, 1).NumberFormat = '@'
$ws.Cells.Item(92, 1).Value = '2026-02-01'
$ws.Cells.Item(92, 2).Value = '21:06:12'
$ws.Cells.Item(92, 3).Value = '21:00'
$ws.Cells.Item(92, 4).Value = 'Bedroom'
$ws.Cells.Item(92, 5).Value = 2
$ws.Cells.Item(92, 6).Value = 'Occupied'
# row 93
$ws.Cells.Item(93, 1).NumberFormat = '@'
$ws.Cells.Item(93, 1).Value = '2026-02-01'
$ws.Cells.Item(93, 2).Value = '21:06:14'
$ws.Cells.Item(93, 3).Value = '21:00'
$ws.Cells.Item(93, 4).Value = 'Bedroom'
$ws.Cells.Item(93, 5).Value = 9
$ws.Cells.Item(93, 6).Value = 'Occupied'
# row 94
$ws.Cells.Item(94, 1).NumberFormat = '@'
$ws.Cells.Item(94, 1).Value = '2026-02-01'
$ws.Cells.Item(94, 2).Value = '21:06:15'
$ws.Cells.Item(94, 3).Value = '21:00'
$ws.Cells.Item(94, 4).Value = 'Bedroom'
$ws.Cells.Item(94, 5).Value = 2
$ws.Cells.Item(94, 6).Value = 'Occupied'
# row 95
$ws.Cells.Item(95, 1).NumberFormat = '@'
$ws.Cells.Item(95, 1).Value = '2026-02-01'
$ws.Cells.Item(95, 2).Value = '21:06:25'
$ws.Cells.Item(95, 3).Value = '21:00'
$ws.Cells.Item(95, 4).Value = 'Bedroom'
$ws.Cells.Item(95, 5).Value = 1
$ws.Cells.Item(95, 6).Value = 'Occupied'

# --- Sheet: mmWave(HR) (Worksheets.Item(9)) ---
$ws = $wb.Worksheets.Item(9)
# row 84
$ws.Cells.Item(84, 1).NumberFormat = '@'
$ws.Cells.Item(84, 1).Value = '2026-02-01'
$ws.Cells.Item(84, 2).Value = '21:05:30'
$ws.Cells.Item(84, 3).Value = '21:00'
$ws.Cells.Item(84, 4).Value = 'Bedroom'
$ws.Cells.Item(84, 5).Value = 53
$ws.Cells.Item(84, 6).Value = 'Occupied'
# row 85
$ws.Cells.Item(85, 1).NumberFormat = '@'
$ws.Cells.Item(85, 1).Value = '2026-02-01'
$ws.Cells.Item(85, 2).Value = '21:05:31'
$ws.Cells.Item(85, 3).Value = '21:00'
$ws.Cells.Item(85, 4).Value = 'Bedroom'
$ws.Cells.Item(85, 5).Value = 50
$ws.Cells.Item(85, 6).Value = 'Occupied'
# row 86
$ws.Cells.Item(86, 1).NumberFormat = '@'
$ws.Cells.Item(86, 1).Value = '2026-02-01'
$ws.Cells.Item(86, 2).Value = '21:05:36'
$ws.Cells.Item(86, 3).Value = '21:00'
$ws.Cells.Item(86, 4).Value = 'Bedroom'
$ws.Cells.Item(86, 5).Value = 49
$ws.Cells.Item(86, 6).Value = 'Occupied'
# row 87
$ws.Cells.Item(87, 1).NumberFormat = '@'
$ws.Cells.Item(87, 1).Value = '2026-02-01'
$ws.Cells.Item(87, 2).Value = '21:05:48'
$ws.Cells.Item(87, 3).Value = '21:00'
$ws.Cells.Item(87, 4).Value = 'Bedroom'
$ws.Cells.Item(87, 5).Value = 50
$ws.Cells.Item(87, 6).Value = 'Occupied'
# row 88
$ws.Cells.Item(88, 1).NumberFormat = '@'
$ws.Cells.Item(88, 1).Value = '2026-02-01'
$ws.Cells.Item(88, 2).Value = '21:05:53'
$ws.Cells.Item(88, 3).Value = '21:00'
$ws.Cells.Item(88, 4).Value = 'Bedroom'
$ws.Cells.Item(88, 5).Value = 49
$ws.Cells.Item(88, 6).Value = 'Occupied'
# row 89
$ws.Cells.Item(89, 1).NumberFormat = '@'
$ws.Cells.Item(89, 1).Value = '2026-02-01'
$ws.Cells.Item(89, 2).Value = '21:06:02'
$ws.Cells.Item(89, 3).Value = '21:00'
$ws.Cells.Item(89, 4).Value = 'Bedroom'
$ws.Cells.Item(89, 5).Value = 67
$ws.Cells.Item(89, 6).Value = 'Occupied'
# row 90
$ws.Cells.Item(90, 1).NumberFormat = '@'
$ws.Cells.Item(90, 1).Value = '2026-02-01'
$ws.Cells.Item(90, 2).Value = '21:06:03'
$ws.Cells.Item(90, 3).Value = '21:00'
$ws.Cells.Item(90, 4).Value = 'Bedroom'
$ws.Cells.Item(90, 5).Value = 51
$ws.Cells.Item(90, 6).Value = 'Occupied'
# row 91
$ws.Cells.Item(91, 1).NumberFormat = '@'
$ws.Cells.Item(91, 1).Value = '2026-02-01'
$ws.Cells.Item(91, 2).Value = '21:06:04'
$ws.Cells.Item(91, 3).Value = '21:00'
$ws.Cells.Item(91, 4).Value = 'Bedroom'
$ws.Cells.Item(91, 5).Value = 50
$ws.Cells.Item(91, 6).Value = 'Occupied'
# row 92
$ws.Cells.Item(92, 1).NumberFormat = '@'
$ws.Cells.Item(92, 1).Value = '2026-02-01'
$ws.Cells.Item(92, 2).Value = '21:06:08'
$ws.Cells.Item(92, 3).Value = '21:00'
$ws.Cells.Item(92, 4).Value = 'Bedroom'
$ws.Cells.Item(92, 5).Value = 49
$ws.Cells.Item(92, 6).Value = 'Occupied'
# row 93
$ws.Cells.Item(93, 1).NumberFormat = '@'
$ws.Cells.Item(93, 1).Value = '2026-02-01'
$ws.Cells.Item(93, 2).Value = '21:06:11'
$ws.Cells.Item(93, 3).Value = '21:00'
$ws.Cells.Item(93, 4).Value = 'Bedroom'
$ws.Cells.Item(93, 5).Value = 50
$ws.Cells.Item(93, 6).Value = 'Occupied'
# row 94
$ws.Cells.Item(94, 1).NumberFormat = '@'
$ws.Cells.Item(94, 1).Value = '2026-02-01'
$ws.Cells.Item(94, 2).Value = '21:06:14'
$ws.Cells.Item(94, 3).Value = '21:00'
$ws.Cells.Item(94, 4).Value = 'Bedroom'
$ws.Cells.Item(94, 5).Value = 57
$ws.Cells.Item(94, 6).Value = 'Occupied'
# row 95
$ws.Cells.Item(95, 1).NumberFormat = '@'
$ws.Cells.Item(95, 1).Value = '2026-02-01'
$ws.Cells.Item(95, 2).Value = '21:06:15'
$ws.Cells.Item(95, 3).Value = '21:00'
$ws.Cells.Item(95, 4).Value = 'Bedroom'
$ws.Cells.Item(95, 5).Value = 50
$ws.Cells.Item(95, 6).Value = 'Occupied'
# row 96
$ws.Cells.Item(96, 1).NumberFormat = '@'
$ws.Cells.Item(96, 1).Value = '2026-02-01'
$ws.Cells.Item(96, 2).Value = '21:06:25'
$ws.Cells.Item(96, 3).Value = '21:00'
$ws.Cells.Item(96, 4).Value = 'Bedroom'
$ws.Cells.Item(96, 5).Value = 49
$ws.Cells.Item(96, 6).Value = 'Occupied'
